$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the data columns (B:Z) to fit the new header row
# (target stored width is 16.7109375; ColumnWidth snaps to the nearest
# pixel-grid step in this host, so 15.83 lands on the closest reachable value)
$ws.Range("B1:Z1").ColumnWidth = 15.83

# Header labels for row 2 (remaining employee/supervisor data structure,
# repeated per year block: supervisor_rating / clients_rating / ai_rating / date)
$headers = @(
    "first_name",
    "id",
    "last_name",
    "age",
    "join_date",
    "supervisor_rating",
    "clients_rating",
    "ai_rating",
    "date",
    "supervisor_rating",
    "clients_rating",
    "ai_rating",
    "date",
    "supervisor_rating",
    "clients_rating",
    "ai_rating",
    "date"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(2, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Font.Size = 11
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $true
}
